# Update countries & provincias Spain
# - Reorders several country rows in the "Pais" sheet so each country name
#   lines up with its correct statistics (the underlying XML reorders the
#   shared-string table; here we express the same end-result by moving the
#   country label together with its row of numbers).
# - Refreshes a couple of headline totals (Estados Unidos / China).
# - Bumps the "Datos actualizados..." timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: the destination Row keeps columns B:H that *belong* to the
# country named `Name`; that data currently lives in row `Src`. Snapshot
# every B:H row we will touch BEFORE writing anything, so reads are never
# affected by earlier writes in this same pass.
$moves = @(
    @{ Row = 144; Name = 'Puerto Rico'; Src = 145 },
    @{ Row = 145; Name = 'Zambia'; Src = 144 },
    @{ Row = 173; Name = 'Granada'; Src = 174 },
    @{ Row = 174; Name = 'Fiyi'; Src = 173 },
    @{ Row = 177; Name = 'Laos'; Src = 178 },
    @{ Row = 178; Name = 'Seychelles'; Src = 177 },
    @{ Row = 179; Name = 'Mozambique'; Src = 180 },
    @{ Row = 180; Name = 'Surinam'; Src = 179 },
    @{ Row = 185; Name = 'San Cristobal y Nieves'; Src = 186 },
    @{ Row = 186; Name = 'Suazilandia'; Src = 185 },
    @{ Row = 187; Name = 'Nepal'; Src = 188 },
    @{ Row = 188; Name = 'Zimbabue'; Src = 187 },
    @{ Row = 193; Name = 'Cabo Verde'; Src = 194 },
    @{ Row = 194; Name = 'San Vicente y las Granadinas'; Src = 193 },
    @{ Row = 200; Name = 'Malaui'; Src = 202 },
    @{ Row = 201; Name = 'Sahara Occidental'; Src = 203 },
    @{ Row = 202; Name = 'Belice'; Src = 200 },
    @{ Row = 203; Name = 'Sierra Leona'; Src = 201 },
    @{ Row = 206; Name = 'Islas Virgenes Britanicas'; Src = 208 },
    @{ Row = 208; Name = 'Anguila'; Src = 206 },
    @{ Row = 210; Name = 'Timor Oriental'; Src = 212 },
    @{ Row = 211; Name = 'Papua Nueva Guinea'; Src = 210 },
    @{ Row = 212; Name = 'Islas Malvinas'; Src = 211 }
)

# Snapshot the B:H values for every distinct source row used above.
$snapshot = @{}
foreach ($m in $moves) {
    $src = [int]$m.Src
    if (-not $snapshot.ContainsKey($src)) {
        $snapshot[$src] = $ws.Range("B$src`:H$src").Value2
    }
}

# Now write the new country name + its matching row of numbers.
foreach ($m in $moves) {
    $row = [int]$m.Row
    $ws.Cells.Item($row, 1).Value = $m.Name
    $ws.Range("B$row`:H$row").Value = $snapshot[[int]$m.Src]
}

# Headline refreshes (not part of the reordering).
$ws.Range("B4").Value = 311357
$ws.Range("C4").Value = 34196
$ws.Range("E4").Value = 288164

$ws.Range("B9").Value = 81669
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = 76964
$ws.Range("E9").Value = 1376
$ws.Range("F9").Value = 295
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 3329

# Timestamp bump.
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 02:52"
